$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Cell value changes -----------------------------------------------
# Order matters for shared-string allocation order, so touch B3 (the
# first "si" to appear) before the two brand new rows.
$ws.Range("B3").Value = "si"

# --- New rule rows ------------------------------------------------------
# Copy formatting from an existing data row so the new rows look the same
# (wrap text, centered/normal style, etc.), then fill in the real values.
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A16:B17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A16").Value = "aplicar_canalizador_localidad"
$ws.Range("A17").Value = "aplicar_canalizador_provincia"
$ws.Range("B16").Value = "si"
$ws.Range("B17").Value = "si"

# Remaining "SI" -> "si" changes
$ws.Range("B14").Value = "si"
$ws.Range("B15").Value = "si"

# --- Column widths --------------------------------------------------
# Column A gets a bit wider and no longer needs the "best fit" flag;
# column B goes back to relying on the sheet's standard width (it no
# longer needs its own explicit width override).
$ws.Columns.Item(1).ColumnWidth = 18.8

# --- Row heights -------------------------------------------------------
# Rows 3, 13, 14 and 15 no longer need the extra-tall 29pt row (their
# text now fits on one line at the new column width), so let Excel
# recompute their height back to the sheet default.
$ws.Rows.Item(3).AutoFit() | Out-Null
$ws.Rows.Item(13).AutoFit() | Out-Null
$ws.Rows.Item(14).AutoFit() | Out-Null
$ws.Rows.Item(15).AutoFit() | Out-Null

# The two new rows wrap onto two lines, same as the other long rules.
$ws.Rows.Item(16).RowHeight = 29
$ws.Rows.Item(17).RowHeight = 29

# --- View / selection ---------------------------------------------------
$ws.Range("A17").Select() | Out-Null
$win = $wb.Windows.Item(1)
$win.ScrollRow = 7
